$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "72.104.35"
Set-TextValue $ws.Range("E2") "  +4.25%  "
Set-TextValue $ws.Range("D3") "4.036.81"
Set-TextValue $ws.Range("E3") "  +3.68%  "
Set-TextValue $ws.Range("E4") "  +0.13%  "
Set-TextValue $ws.Range("D5") "519.65"
Set-TextValue $ws.Range("E5") "  -1.20%  "
Set-TextValue $ws.Range("D6") "147.95"
Set-TextValue $ws.Range("E6") "  +3.08%  "
Set-TextValue $ws.Range("D7") "0.621"
Set-TextValue $ws.Range("E7") "  +2.11%  "
Set-TextValue $ws.Range("E8") "  +0.15%  "
Set-TextValue $ws.Range("D9") "0.734"
Set-TextValue $ws.Range("E9") "  +2.30%  "
Set-TextValue $ws.Range("D10") "0.174"
Set-TextValue $ws.Range("E10") "  +2.40%  "
Set-TextValue $ws.Range("D11") "0.0000333"
Set-TextValue $ws.Range("E11") "  +0.74%  "
Set-TextValue $ws.Range("D12") "47.72"
Set-TextValue $ws.Range("E12") "  +13.85%  "
Set-TextValue $ws.Range("D13") "10.81"
Set-TextValue $ws.Range("E13") "  +6.18%  "
Set-TextValue $ws.Range("D14") "4.679.87"
Set-TextValue $ws.Range("E14") "  +3.69%  "
Set-TextValue $ws.Range("D15") "4.032.26"
Set-TextValue $ws.Range("E15") "  -0.15%  "
Set-TextValue $ws.Range("D16") "21.16"
Set-TextValue $ws.Range("E16") "  +7.69%  "
Set-TextValue $ws.Range("D17") "14.09"
Set-TextValue $ws.Range("E17") "  +2.43%  "
Set-TextValue $ws.Range("E18") "  -1.08%  "
Set-TextValue $ws.Range("E19") "  -2.19%  "
Set-TextValue $ws.Range("D20") "72.052.18"
Set-TextValue $ws.Range("E20") "  +4.29%  "
Set-TextValue $ws.Range("D21") "437.20"
Set-TextValue $ws.Range("E21") "  +3.21%  "
Set-TextValue $ws.Range("D22") "97.36"
Set-TextValue $ws.Range("E22") "  +11.44%  "
Set-TextValue $ws.Range("E23") "  +6.04%  "
Set-TextValue $ws.Range("D24") "14.65"
Set-TextValue $ws.Range("E24") "  +3.92%  "
Set-TextValue $ws.Range("D25") "11.94"
Set-TextValue $ws.Range("E25") "  +3.65%  "
Set-TextValue $ws.Range("D26") "4.01"
Set-TextValue $ws.Range("E26") "  -1.70%  "
Set-TextValue $ws.Range("D27") "11.18"
Set-TextValue $ws.Range("E27") "  +6.03%  "
Set-TextValue $ws.Range("D28") "36.96"
Set-TextValue $ws.Range("E28") "  +3.16%  "
Set-TextValue $ws.Range("D29") "3.10"
Set-TextValue $ws.Range("E29") "  +10.65%  "
Set-TextValue $ws.Range("D30") "696.35"
Set-TextValue $ws.Range("E30") "  +0.10%  "
Set-TextValue $ws.Range("D31") "13.51"
Set-TextValue $ws.Range("E31") "  +3.29%  "
Set-TextValue $ws.Range("D32") "0.129"
Set-TextValue $ws.Range("E32") "  +3.02%  "
Set-TextValue $ws.Range("D33") "7.02"
Set-TextValue $ws.Range("E33") "  +18.12%  "
Set-TextValue $ws.Range("D34") "68.62"
Set-TextValue $ws.Range("E34") "  +1.31%  "
Set-TextValue $ws.Range("D35") "0.0₃0893"
Set-TextValue $ws.Range("E35") "  +6.00%  "
Set-TextValue $ws.Range("D36") "0.437"
Set-TextValue $ws.Range("E36") "  -1.42%  "
Set-TextValue $ws.Range("D37") "3.64"
Set-TextValue $ws.Range("E37") "  +24.48%  "
Set-TextValue $ws.Range("D38") "40.52"
Set-TextValue $ws.Range("E38") "  +1.10%  "
Set-TextValue $ws.Range("E39") "  +4.21%  "
Set-TextValue $ws.Range("D40") "1.00"
Set-TextValue $ws.Range("E40") "  +0.01%  "
Set-TextValue $ws.Range("D41") "0.998"
Set-TextValue $ws.Range("E41") "  -0.11%  "
Set-TextValue $ws.Range("D42") "0.0487"
Set-TextValue $ws.Range("E42") "  +2.00%  "
Set-TextValue $ws.Range("E43") "  +5.50%  "
Set-TextValue $ws.Range("D44") "2.75"
Set-TextValue $ws.Range("E44") "  -1.24%  "
Set-TextValue $ws.Range("E45") "  +4.22%  "
Set-TextValue $ws.Range("D46") "0.145"
Set-TextValue $ws.Range("E46") "  +4.01%  "
Set-TextValue $ws.Range("E47") "  +2.41%  "
Set-TextValue $ws.Range("D48") "9.04"
Set-TextValue $ws.Range("E48") "  +8.74%  "
Set-TextValue $ws.Range("E49") "  +1.48%  "
Set-TextValue $ws.Range("D50") "0.000269"
Set-TextValue $ws.Range("E50") "  +18.78%  "
Set-TextValue $ws.Range("E51") "  +3.60%  "

Write-Output "Done updating cryptos list."
